{"js": "// Locate the paragraph that holds the old text, re-word it, then append a\n// blank paragraph followed by a new paragraph with the reviewable note \u2014\n// exactly mirroring the target diff:\n//   \"To check git diff ingit cli\"  ->  \"This line is edited \"\n//   + <empty paragraph>\n//   + \"New line added to check reviewable\"\n\nconst body = context.document.body;\n\n// Find the run of text that is being edited.\nconst results = body.search(\"To check git diff ingit cli\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Target paragraph text not found\");\n}\n\nconst targetRange = results.items[0];\n\n// Grab a handle on the paragraph that owns the match before the text\n// underneath it changes.\nconst owningParagraphs = targetRange.paragraphs;\nowningParagraphs.load(\"items\");\nawait context.sync();\nconst targetParagraph = owningParagraphs.items[0];\n\n// 1) Replace the paragraph's text.\ntargetRange.insertText(\"This line is edited \", \"Replace\");\nawait context.sync();\n\n// 2) Append a blank paragraph and a new paragraph with the note.\n//    A Flat-OPC insertOoxml() is used (instead of two insertParagraph()\n//    calls) so the blank paragraph is materialised with no <w:r> at all \u2014\n//    matching the byte-for-byte shape produced by the original edit.\nconst flatOpc = `<?xml version=\"1.0\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n<pkg:xmlData>\n<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n<w:body>\n<w:p><w:pPr><w:rPr><w:lang w:val=\"en-US\"/></w:rPr></w:pPr></w:p>\n<w:p><w:pPr><w:rPr><w:lang w:val=\"en-US\"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>New line added to check reviewable</w:t></w:r></w:p>\n</w:body>\n</w:document>\n</pkg:xmlData>\n</pkg:part>\n</pkg:package>`;\n\ntargetParagraph.getRange(\"After\").insertOoxml(flatOpc, \"After\");\nawait context.sync();\n", "ps1": "# Locate the paragraph that holds the old text, re-word it, then append a\n# blank paragraph followed by a new paragraph with the reviewable note \u2014\n# exactly mirroring the target diff:\n#   \"To check git diff ingit cli\"  ->  \"This line is edited \"\n#   + <empty paragraph>\n#   + \"New line added to check reviewable\"\n\n$d = $word.ActiveDocument\n\n# Find the run of text that is being edited; $rng becomes the matched range.\n$rng = $d.Content\n$found = $rng.Find.Execute(\"To check git diff ingit cli\")\nif (-not $found) {\n    throw \"Target paragraph text not found\"\n}\n\n# 1) Replace the paragraph's text.\n$rng.Text = \"This line is edited \"\n\n# 2) Append a blank paragraph and a new paragraph with the note, right after\n#    the edited paragraph. A Flat-OPC InsertXML is used (instead of plain\n#    InsertParagraphAfter calls) so the blank paragraph is materialised with\n#    no run inside it at all \u2014 matching the byte-for-byte shape produced by\n#    the original edit.\n$insertAt = $rng.End\n$target = $d.Range($insertAt, $insertAt)\n\n$flatOpc = @'\n<?xml version=\"1.0\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n<pkg:xmlData>\n<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n<w:body>\n<w:p><w:pPr><w:rPr><w:lang w:val=\"en-US\"/></w:rPr></w:pPr></w:p>\n<w:p><w:pPr><w:rPr><w:lang w:val=\"en-US\"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>New line added to check reviewable</w:t></w:r></w:p>\n</w:body>\n</w:document>\n</pkg:xmlData>\n</pkg:part>\n</pkg:package>\n'@\n\n$target.InsertXML($flatOpc)\n"}
